$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "fricassê"
$ws.Range("E2").Value = "cenoura ralada"

$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "batata e cenoura"

$ws.Range("D4").Value = "kibe"
$ws.Range("E4").Value = "batata doce"

$ws.Range("D5").Value = "frango"
$ws.Range("E5").Value = "abóbora"

$ws.Range("D6").Value = "bife bovino"
$ws.Range("E6").Value = "batata e cenoura"

$ws.Range("B7").Value = "arroz"
$ws.Range("C7").Value = "feijão"
$ws.Range("D7").Value = "ovo cozido"
$ws.Range("E7").Value = "maionese"

$ws.Range("D8").Value = "kibe"
$ws.Range("E8").Value = "batata"
